$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.0059930486527617, 1.058572957734994, 0.30515861163004443, 2.450896347516224)"
$ws.Range("C2").Value = "NIG(1.2863810447940343, 0.9969356232676945, 3.8187708946324035, 6.020804227134572)"
$ws.Range("D2").Value = "NCT(4.131518623008334, 2.432922433847569, -4.18450130471581, 2.907104151795513)"
$ws.Range("E2").Value = "JSU(-1.0570382694231575, 1.1286620342792175, 3.52636161776874, 4.253816246635148)"
